$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column L (12th column) from 26 to 27 (character-width units round
# through a pixel conversion in ColumnWidth, so 26.14 lands on an exact 27)
$ws.Columns.Item(12).ColumnWidth = 26.14

# Refresh the extraction timestamp
$ws.Range("E2").Value = "2026-02-20 11:15:30"

# HUMITAT_MITJANA_DIA: "87%" must remain literal text (same as the source
# inline-string cell), not get auto-converted into a real percentage number.
# Enter it with a quote-prefix to force text, then restore the original
# cell formatting (border/style, no quote-prefix) by copying the format
# from a neighboring plain-text cell in the same row.
$ws.Range("H2").Value = "'87%"
$ws.Range("I2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("J2").Value = "1021.2 hPa"
$ws.Range("K2").Value = "5.0 MJ/m2"
$ws.Range("L2").Value = "11.9 km/h - 249º 10:59 TU"
$ws.Range("M2").Value = "15.6 °C 10:56 TU"
$ws.Range("O2").Value = "4.3 °C"
